$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 475 (shifts existing rows 475.. down by 2)
$ws.Rows.Item(475).Resize(2).Insert(-4121)

# Row 475: new "Primera" record (Packham's Triumph, $/bandeja)
$ws.Range("A475").Value = 10
$ws.Range("B475").Value = "Vega Modelo de Temuco"
$ws.Range("C475").Value = "La Araucanía"
$ws.Range("D475").Value = 44509
$ws.Range("E475").Value = 9
$ws.Range("F475").Value = "Fruta"
$ws.Range("G475").Value = 100104
$ws.Range("H475").Value = "Frutos de pepita"
$ws.Range("I475").Value = 100104005
$ws.Range("J475").Value = "Pera"
$ws.Range("K475").Value = "Packham's Triumph"
$ws.Range("L475").Value = "Primera"
$ws.Range("M475").Value = 125
$ws.Range("N475").Value = 14000
$ws.Range("O475").Value = 14000
$ws.Range("P475").Value = 14000
$ws.Range("Q475").Value = "$/bandeja 18 kilos granel"
$ws.Range("R475").Value = "Región de O'Higgins"
$ws.Range("S475").Value = 778
$ws.Range("T475").Value = 18

# Row 476: new "Segunda" record (Packham's Triumph, $/bandeja)
$ws.Range("A476").Value = 10
$ws.Range("B476").Value = "Vega Modelo de Temuco"
$ws.Range("C476").Value = "La Araucanía"
$ws.Range("D476").Value = 44509
$ws.Range("E476").Value = 9
$ws.Range("F476").Value = "Fruta"
$ws.Range("G476").Value = 100104
$ws.Range("H476").Value = "Frutos de pepita"
$ws.Range("I476").Value = 100104005
$ws.Range("J476").Value = "Pera"
$ws.Range("K476").Value = "Packham's Triumph"
$ws.Range("L476").Value = "Segunda"
$ws.Range("M476").Value = 85
$ws.Range("N476").Value = 12000
$ws.Range("O476").Value = 12000
$ws.Range("P476").Value = 12000
$ws.Range("Q476").Value = "$/bandeja 18 kilos granel"
$ws.Range("R476").Value = "Región de O'Higgins"
$ws.Range("S476").Value = 667
$ws.Range("T476").Value = 18

Write-Output "Done"
